# New wallet-label row: add 2024-10-03 to the Date column (A4).
# Force text entry (matching the existing "Date" column's text-stored
# values) rather than letting Excel auto-convert the literal into a
# date serial, then drop the temporary number-format override so the
# cell ends up unstyled just like its neighbours above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2024-10-03"
$ws.Range("A4").ClearFormats()
